$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Version value (row 3)
$ws.Range("B3").Value = "6.0.0"

# Update Date value (row 8)
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher row (row 9) now has a value
$ws.Range("B9").Value = "Alvearie Team"

# Replace Contact / No display for ContactDetail row (row 10) with Jurisdiction / United States of America
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the duplicate Contact row (old row 11), shifting everything else up
$ws.Rows.Item(11).Delete()
